$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Proportion" column
$ws.Range("H1").Value = "Proportion"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate per-row proportion of portfolio value held in this holding
# (Value of holding / total Value across all holdings for that date)
$ws.Range("H2").Value = 0.2024602832988485
$ws.Range("H3").Value = 0.7975397167011515
$ws.Range("H4").Value = 0.2320683912910165
$ws.Range("H5").Value = 0.7679316087089835
$ws.Range("H6").Value = 0.327064859280383
$ws.Range("H7").Value = 0.672935140719617
$ws.Range("H8").Value = 0.1998387472810739
$ws.Range("H9").Value = 0.8001612527189261
$ws.Range("H10").Value = 0.29448549770074
$ws.Range("H11").Value = 0.70551450229926
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0.285609114310247
$ws.Range("H14").Value = 0.714390885689753
$ws.Range("H15").Value = 0.1899748530660442
$ws.Range("H16").Value = 0.5163135184556231
$ws.Range("H17").Value = 0.2937116284783327
$ws.Range("H18").Value = 0.1965699741749239
$ws.Range("H19").Value = 0.470543036596166
$ws.Range("H20").Value = 0.2395123601042123
$ws.Range("H21").Value = 0.09337462912469786
$ws.Range("H22").Value = 0.2345299855249004
$ws.Range("H23").Value = 0.4117095031671932
$ws.Range("H24").Value = 0.2658091225698614
$ws.Range("H25").Value = 0.08795138873804499
$ws.Range("H26").Value = 0.2202929848295488
$ws.Range("H27").Value = 0.3825462604323004
$ws.Range("H28").Value = 0.3104817565147436
$ws.Range("H29").Value = 0.08667899822340711
$ws.Range("H30").Value = 0.2107788414767302
$ws.Range("H31").Value = 0.3428594152699202
$ws.Range("H32").Value = 0.3735613556684267
$ws.Range("H33").Value = 0.07280038758492284
$ws.Range("H34").Value = 0.3066465430776691
$ws.Range("H35").Value = 0.5792451775289043
$ws.Range("H36").Value = 0.1141082793934266
$ws.Range("H37").Value = 0.2757012999576651
$ws.Range("H38").Value = 0.5538787182061652
$ws.Range("H39").Value = 0.1134214172174498
$ws.Range("H40").Value = 0.05699856461871995
$ws.Range("H41").Value = 0.267700449473472
$ws.Range("H42").Value = 0.4766319116800394
$ws.Range("H43").Value = 0.09544290295554872
$ws.Range("H44").Value = 0.1602247358909399
$ws.Range("H45").Value = 0.1958556412857424
$ws.Range("H46").Value = 0.5775458073398694
$ws.Range("H47").Value = 0.07291873795947541
$ws.Range("H48").Value = 0.1536798134149128
$ws.Range("H49").Value = 0.1606830912658308
$ws.Range("H50").Value = 0.5702535738078549
$ws.Range("H51").Value = 0.06108407087862602
$ws.Range("H52").Value = 0.2079792640476884
$ws.Range("H53").Value = 0.1182621707680098
$ws.Range("H54").Value = 0.6800720369813099
$ws.Range("H55").Value = 0.05303077232766965
$ws.Range("H56").Value = 0.1486350199230106
